# Rename the inline pictures in the document's headers/footers.
#
# The "PearsonLogo.png" pictures (in both footers) are renamed from
# "image2.png" to "image1.png", and the "BTec_Logo-Orange" pictures
# (in both headers) are renamed from "image1.jpg" to "image2.jpg".
#
# Word.Headers/Footers collection indices:
#   Item(1) = wdHeaderFooterPrimary   (default header/footer)
#   Item(2) = wdHeaderFooterFirstPage (first-page header/footer)
#   Item(3) = wdHeaderFooterEvenPages (even-page header/footer, unused here)
#
# Note: InlineShape.Name does not round-trip the name already baked into
# the document (it reads back blank until explicitly (re)assigned), so we
# key off AlternativeText -- which does reflect the picture's description
# -- to find the right picture in each story before renaming it.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
            $shp = $ftr.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -like "*PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
Write-Host "footers done"

for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
            $shp = $hdr.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}
Write-Host "headers done"
